# Updated cryptos list refresh: sync Price (column D) and Volume(1h)
# (column E) cells with the latest scrape. Cells whose new text would
# otherwise be auto-parsed by Excel as a plain number (losing the
# trailing/leading zero formatting of the source text, e.g. "0.940")
# are pinned to a Text number format before the write, then restored to
# the workbook's default ("Normal") style so no stray formatting is left
# behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.385.85'
$ws.Range("E2").Value = '  +0.83%  '
$ws.Range("D3").Value = '1.789.00'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.21'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.89'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.31%  '
$ws.Range("E9").Value = '  +0.86%  '
$ws.Range("E10").Value = '  +0.53%  '
$ws.Range("E11").Value = '  -0.11%  '
$ws.Range("D12").Value = '2.047.18'
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("D14").Value = '1.787.27'
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("E15").Value = '  +1.99%  '
$ws.Range("D16").Value = '34.355.75'
$ws.Range("E16").Value = '  +0.84%  '
$ws.Range("E17").Value = '  +2.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.50'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("E20").Value = '  +1.05%  '
$ws.Range("E21").Value = '  +3.19%  '
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '167.93'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.75%  '
$ws.Range("E25").Value = '  +1.37%  '
$ws.Range("E26").Value = '  +2.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.57'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.84%  '
$ws.Range("E28").Value = '  +1.45%  '
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.04'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.53%  '
$ws.Range("E31").Value = '  +1.86%  '
$ws.Range("E32").Value = '  +2.22%  '
$ws.Range("E33").Value = '  +0.48%  '
$ws.Range("E34").Value = '  +1.70%  '
$ws.Range("E35").Value = '  +6.80%  '
$ws.Range("D36").Value = '1.409.10'
$ws.Range("E36").Value = '  -2.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.683'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.05%  '
$ws.Range("E38").Value = '  +3.12%  '
$ws.Range("E39").Value = '  +0.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '84.18'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.41'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.940'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.65%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.87'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.01%  '
$ws.Range("E45").Value = '  +1.84%  '
$ws.Range("E46").Value = '  +3.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.35%  '
$ws.Range("D48").Value = '1.947.90'
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.39'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.79%  '
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("E51").Value = '  -2.36%  '
